$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so
# values like trailing zeros ("0.0690", "3.60", "1.80") are not
# silently coerced into numbers by Excel's automatic type detection.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.121.07'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.788.54'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.31'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.84'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0690'
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.046.09'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.795.10'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.053.92'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.619'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.03'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.67'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.26'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E26').Value = '  +1.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.29'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.65'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.60'
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.80'
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.452.58'
$ws.Range('E35').Value = '  +4.58%  '
$ws.Range('E36').Value = '  +9.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.645'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('E38').Value = '  +2.51%  '
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.14'
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.53'
$ws.Range('E44').Value = '  +3.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0509'
$ws.Range('E45').Value = '  +2.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.04'
$ws.Range('E46').Value = '  +4.03%  '
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.947.66'
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.17'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('E51').Value = '  -0.04%  '
